$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.718.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.369.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.369.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.946.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.372.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000169"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.855.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.508.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("E28").Value = "  +9.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.770"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.517.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.04%  "
